# Updated resend otp button account block issue
#
# Adds:
#  - sheet1 (loginData): E5 "Incorrrect otp tst", new row 7 (A7 "resend@g.io"
#    as a mailto hyperlink styled like the other account rows, E7 "Resend OTP test")
#  - sheet2 (Message): new row 14 "Max OTP retry limit reached. Please try again later"

$wb = $excel.ActiveWorkbook

# ---- Sheet1 : loginData ----
$ws1 = $wb.Worksheets.Item("loginData")

# New account row 7: resend@g.io (hyperlinked like A2:A6) + its note in E7
$ws1.Range("A7").Value = "resend@g.io"
$ws1.Hyperlinks.Add($ws1.Range("A7"), "mailto:resend@g.io")
$ws1.Range("A7").Style = $ws1.Range("A6").Style

$ws1.Range("E7").Value = "Resend OTP test"

# New "Incorrrect otp tst" note next to row 5 (E5)
$ws1.Range("E5").Value = "Incorrrect otp tst"

# ---- Sheet2 : Message ----
$ws2 = $wb.Worksheets.Item("Message")

# New expected-message row 14
$ws2.Range("A14").Value = "Max OTP retry limit reached. Please try again later"
$ws2.Range("A14").WrapText = $true

# Move sheet1's active-cell selection to E5, then re-select the Message sheet's
# own A14 cell so the workbook keeps "Message" as the active/selected tab
# (matches the original file's tabSelected/activeTab state).
$ws1.Range("E5").Select()
$ws2.Range("A14").Select()
